$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Ben Jarman") {
        $p.Range.Delete()
        break
    }
}
